$d = $word.ActiveDocument

$pairs = @(
    @("681×9=", "471×8="),
    @("171×3=", "388×3="),
    @("924×7=", "638×3="),
    @("981×8=", "429×5="),
    @("115×6=", "324×7="),
    @("301×9=", "474×8="),
    @("518×7=", "671×6="),
    @("333×5=", "262×8="),
    @("245×8=", "324×8="),
    @("464×4=", "117×5="),
    @("368×3=", "278×8="),
    @("582×3=", "580×3="),
    @("404×5=", "128×2="),
    @("313×4=", "788×5="),
    @("577×4=", "784×8="),
    @("709×5=", "526×4="),
    @("374×7=", "804×3="),
    @("450×5=", "275×5="),
    @("585×2=", "486×2="),
    @("608×9=", "430×9="),
    @("863×6=", "843×2="),
    @("863×3=", "171×2="),
    @("393×6=", "843×3="),
    @("858×3=", "418×7="),
    @("992×3=", "626×6=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
